# Insert a new data row at row 158 (pushing the existing rows 158-239 down
# to 159-240) and populate it with the new "Cilantro" price record.
#
# This reproduces the behaviour captured in the diff: a brand new row of
# data appears before what used to be row 158, the dimension grows from
# A1:R239 to A1:R240, and every row from the old 158 through 239 shifts
# down by exactly one row (their content is otherwise untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 158..239 down to 159..240, leaving a blank row 158.
$ws.Rows("158:158").Insert()

# Populate the newly inserted row 158 with the new record's values.
$ws.Range("A158").Value = 11
$ws.Range("B158").Value = "Vega Monumental Concepción"
$ws.Range("C158").Value = "Bíobío"
$ws.Range("D158").Value = 44873
$ws.Range("E158").Value = 8
$ws.Range("F158").Value = 100112040
$ws.Range("G158").Value = "Cilantro"
$ws.Range("H158").Value = "Sin especificar"
$ws.Range("I158").Value = "Primera"
$ws.Range("J158").Value = 130
$ws.Range("K158").Value = 15000
$ws.Range("L158").Value = 16000
$ws.Range("M158").Value = 15385
$ws.Range("N158").Value = "$/caja 36 atados"
$ws.Range("O158").Value = "Región Metropolitana"
$ws.Range("P158").Value = 427
$ws.Range("Q158").Value = 36
$ws.Range("R158").Value = "Hortaliza"
